$d = $word.ActiveDocument

$replacements = @(
    @("117×8=", "647×6="),
    @("853×8=", "221×6="),
    @("527×3=", "562×3="),
    @("972×6=", "400×2="),
    @("135×9=", "260×6="),
    @("347×3=", "514×3="),
    @("465×8=", "927×2="),
    @("510×7=", "296×4="),
    @("405×8=", "728×7="),
    @("916×9=", "978×2="),
    @("107×2=", "173×9="),
    @("977×9=", "880×5="),
    @("830×6=", "435×7="),
    @("584×6=", "294×2="),
    @("951×6=", "493×2="),
    @("354×9=", "624×3="),
    @("935×8=", "530×7="),
    @("401×7=", "478×4="),
    @("667×3=", "638×7="),
    @("654×2=", "617×4="),
    @("473×5=", "649×5="),
    @("373×5=", "961×6="),
    @("215×2=", "374×3="),
    @("487×6=", "135×5="),
    @("359×8=", "863×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
